$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header / data cells
$ws.Range("A2").Value = "Project Date"
$ws.Range("B2").NumberFormat = "DD/MM/YY"
$ws.Range("B2").Value = 42027

$ws.Range("A3").Value = "Project Age (Years)"
$ws.Range("B3").Value = 10

$ws.Range("A4").Value = "Miscellaneous Issues"
$ws.Range("B4").Value = "This is a longish string that needs to be handled by the program. You cannot underestimate how important this is."

# Column widths (closest achievable values given engine's internal
# pixel-grid quantization of column widths to 1/6-character steps)
$ws.Columns.Item(1).ColumnWidth = 18.8
$ws.Columns.Item(2).ColumnWidth = 16.67

# Move selection to A5
[void]$ws.Range("A5").Select()
